{"js": "// Apply the \"Added many more features\" edits to the Jack in a Pot review.\nconst body = context.document.body;\n\n// Each entry is an exact, unique (or intentionally repeated) find/replace\n// pair taken straight from the OOXML diff.\nconst replacements = [\n  {\n    find: \"Play Jack in a Pot for Free - Unique Irish-themed Slot Game\",\n    replace: \"Play Jack in a Pot for Free - Unique Features & Stunning Graphics\",\n  },\n  {\n    find: \"Unique gameplay mechanics and features.\",\n    replace: \"Unique gameplay mechanics\",\n  },\n  {\n    find: \"Visually stunning graphics and animations.\",\n    replace: \"Visually stunning graphics\",\n  },\n  {\n    find: \"Medium volatility with potential for significant payouts.\",\n    replace: \"Medium volatility with potential for big wins\",\n  },\n  {\n    find: \"Themed around Irish folklore and mythology.\",\n    replace: \"Irish-themed with fun symbols and animations\",\n  },\n  {\n    find: \"Low RTP compared to other slot games.\",\n    replace: \"Low theoretical return to player (RTP)\",\n  },\n  {\n    find: \"May not appeal to players who prefer traditional spinning reels.\",\n    replace: \"No spinning reels for traditional slot game enthusiasts\",\n  },\n  {\n    find: \"Discover the unique gameplay mechanics and visually stunning graphics of Jack in a Pot, an Irish-themed slot game. Play for free and win big!\",\n    replace: \"Read our review of Jack in a Pot and play for free! Enjoy unique gameplay mechanics and visually stunning graphics.\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edits to the Jack in a Pot review.\n$d = $word.ActiveDocument\n\n# wdReplaceAll = 2, wdFindContinue = 1 (wrap search through the whole story)\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nfunction Replace-AllText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n}\n\nReplace-AllText \"Play Jack in a Pot for Free - Unique Irish-themed Slot Game\" \"Play Jack in a Pot for Free - Unique Features & Stunning Graphics\"\nReplace-AllText \"Unique gameplay mechanics and features.\" \"Unique gameplay mechanics\"\nReplace-AllText \"Visually stunning graphics and animations.\" \"Visually stunning graphics\"\nReplace-AllText \"Medium volatility with potential for significant payouts.\" \"Medium volatility with potential for big wins\"\nReplace-AllText \"Themed around Irish folklore and mythology.\" \"Irish-themed with fun symbols and animations\"\nReplace-AllText \"Low RTP compared to other slot games.\" \"Low theoretical return to player (RTP)\"\nReplace-AllText \"May not appeal to players who prefer traditional spinning reels.\" \"No spinning reels for traditional slot game enthusiasts\"\nReplace-AllText \"Discover the unique gameplay mechanics and visually stunning graphics of Jack in a Pot, an Irish-themed slot game. Play for free and win big!\" \"Read our review of Jack in a Pot and play for free! Enjoy unique gameplay mechanics and visually stunning graphics.\"\n"}
